$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells with the same style as the existing header row (bold, bordered)
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every data row
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 86
    $ws.Cells.Item($r, 30).Value = 76
    $ws.Cells.Item($r, 31).Value = 0
}
